$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 4.62
$ws.Range("F3").Value = 6.119378603558825
$ws.Range("F4").Value = 6.682026959004283
$ws.Range("F5").Value = 17
$ws.Range("F7").Value = 5.68775654336231
$ws.Range("F8").Value = 3.35
$ws.Range("F9").Value = 1.35705221483896
